# Correctif : mise à jour des moyennes des étudiants (colonne E) suite au
# problème d'insertion de matière en base et à la mise à jour du semestre
# de l'étudiant lors de l'inscription.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPA")

$newAverages = @{
    3  = 11
    4  = 7
    5  = 16
    7  = 17
    8  = 14
    9  = 15
    10 = 6
    11 = 15
    12 = 14
    13 = 16
    14 = 20
    15 = 8
    16 = 17
    17 = 8
    18 = 17
    19 = 18
    21 = 5
    22 = 5
    23 = 8
    24 = 8
    25 = 12
    26 = 11
    27 = 14
    28 = 11
    29 = 11
    30 = 8
    31 = 12
    32 = 11
    33 = 17
    34 = 19
    35 = 15
    36 = 15
    37 = 15
    38 = 11
    39 = 14
    40 = 14
    41 = 18
    42 = 8
    43 = 9
    45 = 11
    46 = 20
    47 = 16
    48 = 8
    49 = 6
    50 = 9
    51 = 7
    52 = 17
    53 = 12
    54 = 9
    55 = 12
    56 = 9
    57 = 11
    58 = 8
    59 = 8
    60 = 7
    61 = 17
    62 = 19
    63 = 17
}

foreach ($row in $newAverages.Keys) {
    $ws.Range("E$row").Value = $newAverages[$row]
}
